$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates ---
$ws.Range("D5").Value = 'Report Generated On: 08/26/2025 09:59 AM'
$ws.Range("C8").Value = 8928.93
$ws.Range("C9").Value = 23
$ws.Range("G10").Value = ""

# --- Remove now-unused trailing rows (old rows 40-51) ---
$ws.Range("A40:I51").EntireRow.Delete()

# --- Rewrite the line-item table (rows 16-38) ---
# Style templates: row 16 (s=9/10/11) for odd positions, row 17 (s=12/13/14) for even positions

$ws.Range("A16:H16").Copy()
$ws.Range("A16:H16").PasteSpecial(-4122)
$ws.Range("A16").Value = 'Point 09'
$ws.Range("B16").Value = 'PLA-HDIG'
$ws.Range("C16").Value = 'Inst'
$ws.Range("D16").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E16").Value = 'EA'
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = 648.53

$ws.Range("A17:H17").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)
$ws.Range("A17").Value = 'Point 17'
$ws.Range("B17").Value = 'GND-CR-4'
$ws.Range("C17").Value = 'Inst'
$ws.Range("D17").Value = 'GND,Cu Clad Rod,#4'
$ws.Range("E17").Value = 'EA'
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = 110.74

$ws.Range("A16:H16").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)
$ws.Range("A18").Value = 'Point 17'
$ws.Range("B18").Value = 'PLA-HDIG'
$ws.Range("C18").Value = 'Inst'
$ws.Range("D18").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E18").Value = 'EA'
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = 648.53

$ws.Range("A17:H17").Copy()
$ws.Range("A19:H19").PasteSpecial(-4122)
$ws.Range("A19").Value = 'Point 19'
$ws.Range("B19").Value = 'PLA-HDIG'
$ws.Range("C19").Value = 'Inst'
$ws.Range("D19").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E19").Value = 'EA'
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = 648.53

$ws.Range("A16:H16").Copy()
$ws.Range("A20:H20").PasteSpecial(-4122)
$ws.Range("A20").Value = 'Point 21'
$ws.Range("B20").Value = 'INS-15-P-S-C'
$ws.Range("C20").Value = 'Inst'
$ws.Range("D20").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("E20").Value = 'EA'
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = 94.17

$ws.Range("A17:H17").Copy()
$ws.Range("A21:H21").PasteSpecial(-4122)
$ws.Range("A21").Value = 'Point 21'
$ws.Range("B21").Value = 'PIN-15-PTP-C'
$ws.Range("C21").Value = 'Inst'
$ws.Range("D21").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("E21").Value = 'EA'
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = 94.17

$ws.Range("A16:H16").Copy()
$ws.Range("A22:H22").PasteSpecial(-4122)
$ws.Range("A22").Value = 'Point 21'
$ws.Range("B22").Value = 'POL-40-2'
$ws.Range("C22").Value = 'Inst'
$ws.Range("D22").Value = 'Pole,40ft,Class 2'
$ws.Range("E22").Value = 'EA'
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = 478.55

$ws.Range("A17:H17").Copy()
$ws.Range("A23:H23").PasteSpecial(-4122)
$ws.Range("A23").Value = 'Point 21'
$ws.Range("B23").Value = 'PLA-HDIG'
$ws.Range("C23").Value = 'Inst'
$ws.Range("D23").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E23").Value = 'EA'
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = 648.53

$ws.Range("A16:H16").Copy()
$ws.Range("A24:H24").PasteSpecial(-4122)
$ws.Range("A24").Value = 'Point 23'
$ws.Range("B24").Value = 'PLA-HDIG'
$ws.Range("C24").Value = 'Inst'
$ws.Range("D24").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E24").Value = 'EA'
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = 648.53

$ws.Range("A17:H17").Copy()
$ws.Range("A25:H25").PasteSpecial(-4122)
$ws.Range("A25").Value = 'Point 25'
$ws.Range("B25").Value = 'GND-CR-4'
$ws.Range("C25").Value = 'Inst'
$ws.Range("D25").Value = 'GND,Cu Clad Rod,#4'
$ws.Range("E25").Value = 'EA'
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = 110.74

$ws.Range("A16:H16").Copy()
$ws.Range("A26:H26").PasteSpecial(-4122)
$ws.Range("A26").Value = 'Point 25'
$ws.Range("B26").Value = 'INS-15-P-S-C'
$ws.Range("C26").Value = 'Inst'
$ws.Range("D26").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("E26").Value = 'EA'
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = ""
$ws.Range("H26").Value = 94.17

$ws.Range("A17:H17").Copy()
$ws.Range("A27:H27").PasteSpecial(-4122)
$ws.Range("A27").Value = 'Point 25'
$ws.Range("B27").Value = 'PIN-15-PTP-C'
$ws.Range("C27").Value = 'Inst'
$ws.Range("D27").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("E27").Value = 'EA'
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = ""
$ws.Range("H27").Value = 94.17

$ws.Range("A16:H16").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)
$ws.Range("A28").Value = 'Point 25'
$ws.Range("B28").Value = 'POL-40-2'
$ws.Range("C28").Value = 'Inst'
$ws.Range("D28").Value = 'Pole,40ft,Class 2'
$ws.Range("E28").Value = 'EA'
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = ""
$ws.Range("H28").Value = 478.55

$ws.Range("A17:H17").Copy()
$ws.Range("A29:H29").PasteSpecial(-4122)
$ws.Range("A29").Value = 'Point 25'
$ws.Range("B29").Value = 'PLA-HDIG'
$ws.Range("C29").Value = 'Inst'
$ws.Range("D29").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E29").Value = 'EA'
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = ""
$ws.Range("H29").Value = 648.53

$ws.Range("A16:H16").Copy()
$ws.Range("A30:H30").PasteSpecial(-4122)
$ws.Range("A30").Value = 'Point 27'
$ws.Range("B30").Value = 'PLA-HDIG'
$ws.Range("C30").Value = 'Inst'
$ws.Range("D30").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E30").Value = 'EA'
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = ""
$ws.Range("H30").Value = 648.53

$ws.Range("A17:H17").Copy()
$ws.Range("A31:H31").PasteSpecial(-4122)
$ws.Range("A31").Value = 'Point 31'
$ws.Range("B31").Value = 'GND-CR-4'
$ws.Range("C31").Value = 'Inst'
$ws.Range("D31").Value = 'GND,Cu Clad Rod,#4'
$ws.Range("E31").Value = 'EA'
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = 110.74

$ws.Range("A16:H16").Copy()
$ws.Range("A32:H32").PasteSpecial(-4122)
$ws.Range("A32").Value = 'Point 31'
$ws.Range("B32").Value = 'PLA-HDIG'
$ws.Range("C32").Value = 'Inst'
$ws.Range("D32").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E32").Value = 'EA'
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = ""
$ws.Range("H32").Value = 648.53

$ws.Range("A17:H17").Copy()
$ws.Range("A33:H33").PasteSpecial(-4122)
$ws.Range("A33").Value = 'Point 33'
$ws.Range("B33").Value = 'PLA-HDIG'
$ws.Range("C33").Value = 'Inst'
$ws.Range("D33").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E33").Value = 'EA'
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = 648.53

$ws.Range("A16:H16").Copy()
$ws.Range("A34:H34").PasteSpecial(-4122)
$ws.Range("A34").Value = 'Point 37'
$ws.Range("B34").Value = 'GND-CR-4'
$ws.Range("C34").Value = 'Inst'
$ws.Range("D34").Value = 'GND,Cu Clad Rod,#4'
$ws.Range("E34").Value = 'EA'
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = ""
$ws.Range("H34").Value = 110.74

$ws.Range("A17:H17").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)
$ws.Range("A35").Value = 'Point 37'
$ws.Range("B35").Value = 'INS-15-P-S-C'
$ws.Range("C35").Value = 'Inst'
$ws.Range("D35").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("E35").Value = 'EA'
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = ""
$ws.Range("H35").Value = 94.17

$ws.Range("A16:H16").Copy()
$ws.Range("A36:H36").PasteSpecial(-4122)
$ws.Range("A36").Value = 'Point 37'
$ws.Range("B36").Value = 'PIN-15-PTP-C'
$ws.Range("C36").Value = 'Inst'
$ws.Range("D36").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("E36").Value = 'EA'
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = ""
$ws.Range("H36").Value = 94.17

$ws.Range("A17:H17").Copy()
$ws.Range("A37:H37").PasteSpecial(-4122)
$ws.Range("A37").Value = 'Point 37'
$ws.Range("B37").Value = 'POL-40-2'
$ws.Range("C37").Value = 'Inst'
$ws.Range("D37").Value = 'Pole,40ft,Class 2'
$ws.Range("E37").Value = 'EA'
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = ""
$ws.Range("H37").Value = 478.55

$ws.Range("A16:H16").Copy()
$ws.Range("A38:H38").PasteSpecial(-4122)
$ws.Range("A38").Value = 'Point 37'
$ws.Range("B38").Value = 'PLA-HDIG'
$ws.Range("C38").Value = 'Inst'
$ws.Range("D38").Value = 'PLA,Hand Dig or Additional  Excavation'
$ws.Range("E38").Value = 'EA'
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = ""
$ws.Range("H38").Value = 648.53

# --- TOTAL row ---
$ws.Range("A39").Value = "TOTAL"
$ws.Range("H39").Value = 8928.929999999998
$ws.Range("A39:G39").Merge()
$excel.Application.CutCopyMode = $false
